$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Robot22 -> Robot32, (6, 1) -> (2, 9)
$ws.Range("A1").Value = "Move Robot32 to location (2, 9) and remove the toolkit."
$ws.Range("B1").Value = "['Robot32']"
$ws.Range("E1").Value = "(2, 9)"

# Row 2: Robot6 stays, (3, 8) -> (5, 4)
$ws.Range("A2").Value = "Move Robot6 to location (5, 4) and remove the liquid spill."
$ws.Range("E2").Value = "(5, 4)"

# Row 3: Robot32 -> Robot29, gripper -> front loader, (12, 1) -> (5, 12)
$ws.Range("A3").Value = "Move Robot29 to location (5, 12) and remove the large debris."
$ws.Range("B3").Value = "['Robot29']"
$ws.Range("C3").Value = "['front loader']"
$ws.Range("E3").Value = "(5, 12)"

# Row 4: Robot28 -> Robot48, (9, 4) -> (6, 6)
$ws.Range("A4").Value = "Move Robot48 to location (6, 6) and remove the dust."
$ws.Range("B4").Value = "['Robot48']"
$ws.Range("E4").Value = "(6, 6)"

# Row 5: Robot41 stays, (7, 6) -> (1, 8)
$ws.Range("A5").Value = "Move Robot41 to location (1, 8) and remove the grass."
$ws.Range("E5").Value = "(1, 8)"

# Row 6: Robot50 -> Robot10, (12, 12) -> (9, 5)
$ws.Range("A6").Value = "Move Robot10 to location (9, 5) and remove the small debris."
$ws.Range("B6").Value = "['Robot10']"
$ws.Range("E6").Value = "(9, 5)"

# Row 7: Robot23 -> Robot13, (9, 4) -> (10, 10)
$ws.Range("A7").Value = "Move Robot13 to location (10, 10) and remove the vehicle."
$ws.Range("B7").Value = "['Robot13']"
$ws.Range("E7").Value = "(10, 10)"

# Row 8: Robot42 -> Robot23, gripper -> tow hook, (8, 8) -> (8, 2)
$ws.Range("A8").Value = "Move Robot23 to location (8, 2) and remove the construction materials."
$ws.Range("B8").Value = "['Robot23']"
$ws.Range("C8").Value = "['tow hook']"
$ws.Range("E8").Value = "(8, 2)"

# Row 9: Robot24 stays, (6, 10) -> (11, 8)
$ws.Range("A9").Value = "Move Robot24 to location (11, 8) and remove the tree branches."
$ws.Range("E9").Value = "(11, 8)"

# Row 10: Robot15 stays, (9, 4) -> (3, 8)
$ws.Range("A10").Value = "Move Robot15 to location (3, 8) and remove the screws."
$ws.Range("E10").Value = "(3, 8)"
